# [ADD] Extended map + new url for connection
#
# Fills in the new X_coordinate (F) / Y_coordinate (G) columns for rows
# 34-47 of the map, and marks the single outlier cell (G38) with the
# custom "#,##0 _€;-#,##0 _€" number format that the author introduced
# for this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [X_coordinate (F), Y_coordinate (G)]
$rows = @{
    34 = @(159, 470)
    35 = @(214, 439)
    36 = @(104, 441)
    37 = @(714, 563)
    38 = @(657, 471)
    39 = @(602, 407)
    40 = @(521, 462)
    41 = @(462, 517)
    42 = @(517, 607)
    43 = @(575, 611)
    44 = @(659, 627)
    45 = @(631, 548)
    46 = @(599, 503)
    47 = @(549, 536)
}

foreach ($r in $rows.Keys) {
    $xy = $rows[$r]
    $ws.Cells.Item($r, 6).Value = $xy[0]
    $ws.Cells.Item($r, 7).Value = $xy[1]
}

# G38 carries the new custom currency-ish number format (numFmtId 165)
$ws.Cells.Item(38, 7).NumberFormat = "#,##0\ _€;\-#,##0\ _€"

# Restore the author's final cursor position / scroll state.
$ws.Range("G47").Select()

$wb.Save()
